$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '42.652.82'
$ws.Range("E2").Value = '  -1.34%  '
$ws.Range("D3").Value = '2.375.42'
$ws.Range("E3").Value = '  +0.20%  '
$ws.Range("E4").Value = '  -0.18%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '330.70'
$ws.Range("E5").Value = '  +6.31%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '99.65'
$ws.Range("E6").Value = '  -7.29%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.637'
$ws.Range("E7").Value = '  -0.64%  '
$ws.Range("E8").Value = '  +0.01%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.622'
$ws.Range("E9").Value = '  -1.60%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '40.22'
$ws.Range("E10").Value = '  -6.93%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0920'
$ws.Range("E11").Value = '  -2.20%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '8.49'
$ws.Range("E12").Value = '  -5.75%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '1.01'
$ws.Range("E13").Value = '  -5.83%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.106'
$ws.Range("E14").Value = '  +0.44%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '16.35'
$ws.Range("E15").Value = '  -1.32%  '
$ws.Range("D16").Value = '2.719.82'
$ws.Range("E16").Value = '  -0.17%  '
$ws.Range("D17").Value = '2.366.58'
$ws.Range("E17").Value = '  +0.40%  '
$ws.Range("D18").Value = '42.608.38'
$ws.Range("E18").Value = '  -1.47%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '7.82'
$ws.Range("E19").Value = '  +5.80%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.0000107'
$ws.Range("E20").Value = '  -2.43%  '
$ws.Range("E21").Value = '  +7.97%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '75.08'
$ws.Range("E22").Value = '  -0.65%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '275.91'
$ws.Range("E23").Value = '  +8.92%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.30'
$ws.Range("E24").Value = '  -8.35%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '9.69'
$ws.Range("E25").Value = '  +7.80%  '
$ws.Range("E26").Value = '  +0.09%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '11.47'
$ws.Range("B28").Value = 'LEO'
$ws.Range("C28").Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '3.97'
$ws.Range("E28").Value = '  -0.55%  '
$ws.Range("B29").Value = 'EthereumClassic'
$ws.Range("C29").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '23.80'
$ws.Range("E29").Value = '  +4.24%  '
$ws.Range("B30").Value = 'Toncoin'
$ws.Range("C30").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '2.21'
$ws.Range("E30").Value = '  -2.86%  '
$ws.Range("B31").Value = 'Monero'
$ws.Range("C31").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '174.89'
$ws.Range("E31").Value = '  +0.97%  '
$ws.Range("B32").Value = 'WEMIXToken'
$ws.Range("C32").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.10'
$ws.Range("E32").Value = '  -2.33%  '
$ws.Range("B33").Value = 'Hedera'
$ws.Range("C33").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.0903'
$ws.Range("E33").Value = '  -0.82%  '
$ws.Range("B34").Value = 'InjectiveProtocol'
$ws.Range("C34").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '35.37'
$ws.Range("E34").Value = '  -10.21%  '
$ws.Range("B35").Value = 'Filecoin'
$ws.Range("C35").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '6.05'
$ws.Range("E35").Value = '  +3.17%  '
$ws.Range("B36").Value = 'Stellar'
$ws.Range("C36").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.134'
$ws.Range("E36").Value = '  +1.59%  '
$ws.Range("B37").Value = 'RenderToken'
$ws.Range("C37").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '4.57'
$ws.Range("E37").Value = '  -8.67%  '
$ws.Range("B38").Value = 'VeChain'
$ws.Range("C38").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.0359'
$ws.Range("E38").Value = '  -5.09%  '
$ws.Range("B39").Value = 'LidoDAOToken'
$ws.Range("C39").Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '2.96'
$ws.Range("E39").Value = '  +6.62%  '
$ws.Range("B40").Value = 'NEARProtocol'
$ws.Range("C40").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '3.87'
$ws.Range("E40").Value = '  -5.47%  '
$ws.Range("B41").Value = 'Kaspa'
$ws.Range("C41").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.106'
$ws.Range("E41").Value = '  +1.44%  '
$ws.Range("B42").Value = 'ARBITRUM'
$ws.Range("C42").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.53'
$ws.Range("E42").Value = '  -1.02%  '
$ws.Range("B43").Value = 'MultiversX'
$ws.Range("C43").Value = 'https://coinranking.com/coin/omwkOTglq+multiversx-egld'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '69.63'
$ws.Range("E43").Value = '  -4.12%  '
$ws.Range("B44").Value = 'Algorand'
$ws.Range("C44").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.229'
$ws.Range("E44").Value = '  -2.17%  '
$ws.Range("B45").Value = 'FirstDigitalUSD'
$ws.Range("C45").Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '1.00'
$ws.Range("E45").Value = '  -0.19%  '
$ws.Range("B46").Value = 'BitcoinSV'
$ws.Range("C46").Value = 'https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '92.47'
$ws.Range("E46").Value = '  +33.41%  '
$ws.Range("B47").Value = 'Aave'
$ws.Range("C47").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '116.74'
$ws.Range("E47").Value = '  +3.58%  '
$ws.Range("B48").Value = 'Celestia'
$ws.Range("C48").Value = 'https://coinranking.com/coin/YQcD0lBl7+celestia-tia'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '11.99'
$ws.Range("E48").Value = '  -3.69%  '
$ws.Range("B49").Value = 'THORChain'
$ws.Range("C49").Value = 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '5.48'
$ws.Range("E49").Value = '  -3.72%  '
$ws.Range("B50").Value = 'FraxShare'
$ws.Range("C50").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '9.11'
$ws.Range("E50").Value = '  -3.07%  '
$ws.Range("B51").Value = 'Maker'
$ws.Range("C51").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D51").Value = '1.599.18'
$ws.Range("E51").Value = '  +6.75%  '
